$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "(203957296, Omri Ben Shabat: -4,-7)"
$ws.Range("B1").Value = "(206532695, Matan Vakrat: -1,-8)"
$ws.Range("C1").Value = "(302962915, Asher  Odeh: -9,-9)"
$ws.Range("D1").Value = "(308035542, Anastasia  Kubi: -10,-2)"
$ws.Range("E1").Value = "(311177802, Christina  Uksusman: 8,4)"
$ws.Range("F1").Value = "(305251175, Or  Leder: -3,-6)"

$ws.Range("A3").Value = "cost: 524.7044400762932"
$ws.Range("A4").Value = "time: 62.46305500953665"
